$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to match repulled data
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = 0
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = -4
$ws.Range("F13").Value = -7
$ws.Range("F14").Value = -1
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = -3
